$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F6: date changed from 2024-05-05 to 2026-03-05 (kept as text, not an Excel date)
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "2026-03-05"
$ws.Range("F6").ClearFormats()

# Row 8 is edited in place to hold the data that used to live in row 10
# (id 9 / ivan / gorda / 460036), while its Fecha becomes the old row 10's
# date (2024-12-11). E8 (Tarjeta) is untouched - it already matches.
$ws.Range("A8").Value = 9
$ws.Range("B8").Value = "ivan"
$ws.Range("C8").Value = "gorda"
$ws.Range("D8").Value = 460036

$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "2024-12-11"
$ws.Range("F8").ClearFormats()

# Rows 9 and 10 (the old "maxi morales" and "ivan gorda" records) are
# removed entirely - row 10's data now lives in row 8 above.
$ws.Range("A9:F10").Delete()
